$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scenario text in A3: the booking destination changed
# from Mayiladuthurai to Karaikal.
$ws.Range("A3").Value = "From Chennai to Karaikal"

# Update the selected/active cell shown in the sheet view.
$ws.Range("C14").Select()
